$d = $word.ActiveDocument

$replacements = @(
    @{Old = "2024-08-19 Monday"; New = "2024-08-20 Tuesday"},
    @{Old = "939÷6="; New = "232÷2="},
    @{Old = "974÷6="; New = "435÷3="},
    @{Old = "140÷8="; New = "920÷9="},
    @{Old = "769÷5="; New = "206÷4="},
    @{Old = "266÷6="; New = "987÷7="},
    @{Old = "681÷6="; New = "149÷5="},
    @{Old = "933÷9="; New = "637÷6="},
    @{Old = "574÷6="; New = "125÷8="},
    @{Old = "102÷3="; New = "857÷5="},
    @{Old = "787÷6="; New = "461÷9="},
    @{Old = "376÷4="; New = "665÷6="},
    @{Old = "454÷4="; New = "251÷3="},
    @{Old = "250÷6="; New = "436÷4="},
    @{Old = "517÷2="; New = "232÷4="},
    @{Old = "908÷4="; New = "398÷6="},
    @{Old = "586÷2="; New = "323÷3="},
    @{Old = "627÷2="; New = "887÷9="},
    @{Old = "310÷4="; New = "655÷9="},
    @{Old = "766÷2="; New = "122÷3="},
    @{Old = "714÷9="; New = "956÷8="},
    @{Old = "432÷6="; New = "367÷6="},
    @{Old = "825÷5="; New = "456÷9="},
    @{Old = "646÷2="; New = "380÷5="},
    @{Old = "436÷3="; New = "963÷4="},
    @{Old = "436÷5="; New = "671÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
